# Issue 9 Status: Verified — Documentación de las nuevas interfaces y uso de pines
# Adds the new "Solenoide" 3V commands (activar/desactivar/comprobar) to the
# "Ordenes" sheet (rows 17-19, reusing the reserved 0x67/0x68/0x69 codes) and
# labels the two newly-used digital pins (D10/D11) on the "Pines" sheet.

$wb  = $excel.ActiveWorkbook
$ordenes = $wb.Worksheets.Item("Ordenes")
$pines   = $wb.Worksheets.Item("Pines")

# ---------------------------------------------------------------------------
# "Ordenes" sheet — rows 17/18/19 currently only hold the reserved command
# letter in column H (g/h/i); fill in the command name, hex code, decimal
# code and the G-column mirror of the letter, matching rows 11-16 above them.
# ---------------------------------------------------------------------------

# Row 17 — activarSolenoide3V() / 0x67 / 103 / g
$ordenes.Range("A17").Value = "activarSolenoide3V()"
$ordenes.Range("A17").Font.Color = 255
$ordenes.Range("E17").Value = "0x67"
$ordenes.Range("F17").Value = 103

$ordenes.Range("G14").Copy() | Out-Null
$ordenes.Range("G17").PasteSpecial(-4122) | Out-Null
$ordenes.Range("G17").Value = "g"

$ordenes.Range("H11").Copy() | Out-Null
$ordenes.Range("H17").PasteSpecial(-4122) | Out-Null

# Row 18 — desactivarSolenoide3V() / 0x68 / 104 / h
$ordenes.Range("A18").Value = "desactivarSolenoide3V()"
$ordenes.Range("A18").Font.Color = 255
$ordenes.Range("E18").Value = "0x68"
$ordenes.Range("F18").Value = 104

$ordenes.Range("G14").Copy() | Out-Null
$ordenes.Range("G18").PasteSpecial(-4122) | Out-Null
$ordenes.Range("G18").Value = "h"

$ordenes.Range("H11").Copy() | Out-Null
$ordenes.Range("H18").PasteSpecial(-4122) | Out-Null

# Row 19 — comprobarSolenoide3V() / 0x69 / 105 / i
$ordenes.Range("A19").Value = "comprobarSolenoide3V()"
$ordenes.Range("A19").Font.Color = 255
$ordenes.Range("E19").Value = "0x69"
$ordenes.Range("F19").Value = 105

$ordenes.Range("G14").Copy() | Out-Null
$ordenes.Range("G19").PasteSpecial(-4122) | Out-Null
$ordenes.Range("G19").Value = "i"

$ordenes.Range("H11").Copy() | Out-Null
$ordenes.Range("H19").PasteSpecial(-4122) | Out-Null

$excel.CutCopyMode = $false

# ---------------------------------------------------------------------------
# "Pines" sheet — label the pins (D10/D11) now used for the Solenoide signal.
# ---------------------------------------------------------------------------

$pines.Range("D14").Copy() | Out-Null
$pines.Range("D19").PasteSpecial(-4122) | Out-Null
$pines.Range("C19").Value = "SolenoideOn"

$pines.Range("D14").Copy() | Out-Null
$pines.Range("D20").PasteSpecial(-4122) | Out-Null
$pines.Range("C20").Value = "SolenoideOff"

$excel.CutCopyMode = $false

# ---------------------------------------------------------------------------
# Restore the cursor / selection on each sheet (cosmetic, matches the saved
# view state) — select the inactive sheet first so the active-tab stays put.
# ---------------------------------------------------------------------------

$pines.Range("F20").Select() | Out-Null
$ordenes.Range("P16").Select() | Out-Null
